$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.574.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.566.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.493'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +5.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0885'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.790.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.564.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.52%  '
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.570.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  -5.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("E29").Value = '  -3.57%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.398.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.69%  '
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0166'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.536'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.790'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("E44").Value = '  +2.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.976'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.702.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("E50").Value = '  -4.64%  '
$ws.Range("E51").Value = '  -0.88%  '
